# "Fichiers pour la map"
# The original deck has five straight-line connectors (the map's grid
# lines) whose outline weight is 4.5pt (57150 EMU). Thin them down to
# 0.25pt (3175 EMU), matching the canonical diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$msoConnector = 9

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Type -eq $msoConnector -and [math]::Abs($sh.Line.Weight - 4.5) -lt 0.01) {
        $sh.Line.Weight = 0.25
    }
}
